$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "22.077.11"
$ws.Range("E2").Value = "  -1.41%  "

# Row 3
$ws.Range("D3").Value = "1.557.83"
$ws.Range("E3").Value = "  -0.48%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9999"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.55%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3867"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3244"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.86%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.09%  "

# Row 10
$ws.Range("E10").Value = "  -1.88%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07367"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.53%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9999"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.94%  "

# Row 14
$ws.Range("E14").Value = "  -2.32%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.805"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.18%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001126"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.83%  "

# Row 17
$ws.Range("D17").Value = "1.554.69"
$ws.Range("E17").Value = "  -0.67%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06609"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.20%  "

# Row 19
$ws.Range("E19").Value = "  -0.84%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.406"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.27%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9993"
$ws.Range("D21").Style = "Normal"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.36%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.34%  "

# Row 24
$ws.Range("D24").Value = "22.088.53"
$ws.Range("E24").Value = "  -1.32%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.341"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.38%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.553"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.04%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.02%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.865"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.46%  "

# Row 30
$ws.Range("D30").Value = "1.730.64"
$ws.Range("E30").Value = "  -0.60%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.14%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.109"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.63%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.836"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.60%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.702"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -12.85%  "

# Row 35
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.359"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.69%  "

# Row 36
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08210"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06278"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.06%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02304"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.87%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.233"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.17%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2114"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.06%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.223"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.36%  "

# Row 42
$ws.Range("E42").Value = "  -1.93%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9994"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.11%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5960"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.26%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.94%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.713"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.80%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5757"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.06%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.933"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.47%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "119.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.75%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.157"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.59%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06891"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.62%  "

